# Auto-generated Excel COM-interop script applying the diff
# to the Pandaemonium Profits workbook (sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# ALC row 11
$ws_ALC.Cells.Item(11, 8).Value = 1255.2106
$ws_ALC.Cells.Item(11, 9).Value = 1255.2106
$ws_ALC.Cells.Item(11, 10).Value = 0
$ws_ALC.Cells.Item(11, 11).Value = 1255.2106
$ws_ALC.Cells.Item(11, 12).Value = 0
$ws_ALC.Cells.Item(11, 13).Value = -1115.2106

# ALC row 40
$ws_ALC.Cells.Item(40, 8).Value = 3009
$ws_ALC.Cells.Item(40, 9).Value = 4975
$ws_ALC.Cells.Item(40, 10).Value = 1698.3334
$ws_ALC.Cells.Item(40, 11).Value = 4975
$ws_ALC.Cells.Item(40, 12).Value = 1698.3334
$ws_ALC.Cells.Item(40, 13).Value = -4800
$ws_ALC.Cells.Item(40, 14).Value = -2048.3334

# ALC row 46
$ws_ALC.Cells.Item(46, 8).Value = 10102800
$ws_ALC.Cells.Item(46, 9).Value = 3533.3333
$ws_ALC.Cells.Item(46, 10).Value = 14431057
$ws_ALC.Cells.Item(46, 11).Value = 10599.9999
$ws_ALC.Cells.Item(46, 12).Value = 43293171
$ws_ALC.Cells.Item(46, 13).Value = -10480.9999
$ws_ALC.Cells.Item(46, 14).Value = -43293409

# ALC row 60
$ws_ALC.Cells.Item(60, 8).Value = 10102800
$ws_ALC.Cells.Item(60, 9).Value = 3533.3333
$ws_ALC.Cells.Item(60, 10).Value = 14431057
$ws_ALC.Cells.Item(60, 11).Value = 10599.9999
$ws_ALC.Cells.Item(60, 12).Value = 43293171
$ws_ALC.Cells.Item(60, 13).Value = -10115.9999
$ws_ALC.Cells.Item(60, 14).Value = -43294139

# ALC row 69
$ws_ALC.Cells.Item(69, 8).Value = 76555.57000000001
$ws_ALC.Cells.Item(69, 9).Value = 6009.2856
$ws_ALC.Cells.Item(69, 10).Value = 147101.86
$ws_ALC.Cells.Item(69, 11).Value = 18027.8568
$ws_ALC.Cells.Item(69, 12).Value = 441305.58
$ws_ALC.Cells.Item(69, 13).Value = -17153.8568
$ws_ALC.Cells.Item(69, 14).Value = -443053.58

# ALC row 72
$ws_ALC.Cells.Item(72, 8).Value = 76555.57000000001
$ws_ALC.Cells.Item(72, 9).Value = 6009.2856
$ws_ALC.Cells.Item(72, 10).Value = 147101.86
$ws_ALC.Cells.Item(72, 11).Value = 54083.5704
$ws_ALC.Cells.Item(72, 12).Value = 1323916.74
$ws_ALC.Cells.Item(72, 13).Value = -49715.5704
$ws_ALC.Cells.Item(72, 14).Value = -1332652.74

# ALC row 135
$ws_ALC.Cells.Item(135, 8).Value = 55557880
$ws_ALC.Cells.Item(135, 9).Value = 23811466
$ws_ALC.Cells.Item(135, 10).Value = 166670320
$ws_ALC.Cells.Item(135, 11).Value = 214303194
$ws_ALC.Cells.Item(135, 12).Value = 1500032880
$ws_ALC.Cells.Item(135, 13).Value = -214300659
$ws_ALC.Cells.Item(135, 14).Value = -1500037950

# ALC row 138
$ws_ALC.Cells.Item(138, 8).Value = 3287.99
$ws_ALC.Cells.Item(138, 9).Value = 1791.2963
$ws_ALC.Cells.Item(138, 10).Value = 3841.5615
$ws_ALC.Cells.Item(138, 11).Value = 5373.8889
$ws_ALC.Cells.Item(138, 12).Value = 11524.6845
$ws_ALC.Cells.Item(138, 13).Value = -233.8888999999999
$ws_ALC.Cells.Item(138, 14).Value = -21804.6845

# ALC row 141
$ws_ALC.Cells.Item(141, 8).Value = 2134.7917
$ws_ALC.Cells.Item(141, 9).Value = 1532.0588
$ws_ALC.Cells.Item(141, 10).Value = 3598.5715
$ws_ALC.Cells.Item(141, 11).Value = 4596.1764
$ws_ALC.Cells.Item(141, 12).Value = 10795.7145
$ws_ALC.Cells.Item(141, 13).Value = 583.8235999999997
$ws_ALC.Cells.Item(141, 14).Value = -21155.7145

# ARM row 24
$ws_ARM.Cells.Item(24, 8).Value = 24813.5
$ws_ARM.Cells.Item(24, 9).Value = 0
$ws_ARM.Cells.Item(24, 10).Value = 24813.5
$ws_ARM.Cells.Item(24, 11).Value = 0
$ws_ARM.Cells.Item(24, 12).Value = 24813.5
$ws_ARM.Cells.Item(24, 14).Value = -25561.5

# ARM row 97
$ws_ARM.Cells.Item(97, 8).Value = 982.8570999999999
$ws_ARM.Cells.Item(97, 9).Value = 824.36
$ws_ARM.Cells.Item(97, 10).Value = 2303.6667
$ws_ARM.Cells.Item(97, 11).Value = 824.36
$ws_ARM.Cells.Item(97, 12).Value = 2303.6667
$ws_ARM.Cells.Item(97, 13).Value = -328.36
$ws_ARM.Cells.Item(97, 14).Value = -3295.6667

# ARM row 100
$ws_ARM.Cells.Item(100, 8).Value = 24813.5
$ws_ARM.Cells.Item(100, 9).Value = 0
$ws_ARM.Cells.Item(100, 10).Value = 24813.5
$ws_ARM.Cells.Item(100, 11).Value = 0
$ws_ARM.Cells.Item(100, 12).Value = 24813.5
$ws_ARM.Cells.Item(100, 14).Value = -26977.5

# BSM row 99
$ws_BSM.Cells.Item(99, 8).Value = 2065.7144
$ws_BSM.Cells.Item(99, 9).Value = 1845.8823
$ws_BSM.Cells.Item(99, 10).Value = 3000
$ws_BSM.Cells.Item(99, 11).Value = 1845.8823
$ws_BSM.Cells.Item(99, 12).Value = 3000
$ws_BSM.Cells.Item(99, 13).Value = -347.8823
$ws_BSM.Cells.Item(99, 14).Value = -5996

# CRP row 4
$ws_CRP.Cells.Item(4, 8).Value = 9664.200000000001
$ws_CRP.Cells.Item(4, 9).Value = 0
$ws_CRP.Cells.Item(4, 10).Value = 9664.200000000001
$ws_CRP.Cells.Item(4, 11).Value = 0
$ws_CRP.Cells.Item(4, 12).Value = 9664.200000000001
$ws_CRP.Cells.Item(4, 14).Value = -9888.200000000001

# CRP row 31
$ws_CRP.Cells.Item(31, 8).Value = 3163.8235
$ws_CRP.Cells.Item(31, 9).Value = 0
$ws_CRP.Cells.Item(31, 10).Value = 3163.8235
$ws_CRP.Cells.Item(31, 11).Value = 0
$ws_CRP.Cells.Item(31, 12).Value = 3163.8235
$ws_CRP.Cells.Item(31, 13).ClearContents()
$ws_CRP.Cells.Item(31, 14).Value = -3753.8235

# CRP row 34
$ws_CRP.Cells.Item(34, 8).Value = 3163.8235
$ws_CRP.Cells.Item(34, 9).Value = 0
$ws_CRP.Cells.Item(34, 10).Value = 3163.8235
$ws_CRP.Cells.Item(34, 11).Value = 0
$ws_CRP.Cells.Item(34, 12).Value = 3163.8235
$ws_CRP.Cells.Item(34, 13).ClearContents()
$ws_CRP.Cells.Item(34, 14).Value = -3567.8235

# CRP row 53
$ws_CRP.Cells.Item(53, 8).Value = 47254.75
$ws_CRP.Cells.Item(53, 9).Value = 0
$ws_CRP.Cells.Item(53, 10).Value = 47254.75
$ws_CRP.Cells.Item(53, 11).Value = 0
$ws_CRP.Cells.Item(53, 12).Value = 47254.75
$ws_CRP.Cells.Item(53, 14).Value = -48468.75

# CRP row 60
$ws_CRP.Cells.Item(60, 8).Value = 16459.8
$ws_CRP.Cells.Item(60, 9).Value = 12933
$ws_CRP.Cells.Item(60, 10).Value = 21750
$ws_CRP.Cells.Item(60, 11).Value = 12933
$ws_CRP.Cells.Item(60, 12).Value = 21750
$ws_CRP.Cells.Item(60, 13).Value = -12422
$ws_CRP.Cells.Item(60, 14).Value = -22772

# CRP row 86
$ws_CRP.Cells.Item(86, 8).Value = 2440.8572
$ws_CRP.Cells.Item(86, 9).Value = 1696.5
$ws_CRP.Cells.Item(86, 10).Value = 3433.3333
$ws_CRP.Cells.Item(86, 11).Value = 1696.5
$ws_CRP.Cells.Item(86, 12).Value = 3433.3333
$ws_CRP.Cells.Item(86, 13).Value = -573.5
$ws_CRP.Cells.Item(86, 14).Value = -5679.3333

# CRP row 89
$ws_CRP.Cells.Item(89, 8).Value = 2440.8572
$ws_CRP.Cells.Item(89, 9).Value = 1696.5
$ws_CRP.Cells.Item(89, 10).Value = 3433.3333
$ws_CRP.Cells.Item(89, 11).Value = 8482.5
$ws_CRP.Cells.Item(89, 12).Value = 17166.6665
$ws_CRP.Cells.Item(89, 13).Value = -2866.5
$ws_CRP.Cells.Item(89, 14).Value = -28398.6665

# CUL row 26
$ws_CUL.Cells.Item(26, 8).Value = 575.9231
$ws_CUL.Cells.Item(26, 9).Value = 83
$ws_CUL.Cells.Item(26, 10).Value = 998.4286
$ws_CUL.Cells.Item(26, 11).Value = 249
$ws_CUL.Cells.Item(26, 12).Value = 2995.2858
$ws_CUL.Cells.Item(26, 13).Value = 39
$ws_CUL.Cells.Item(26, 14).Value = -3571.2858

# CUL row 74
$ws_CUL.Cells.Item(74, 8).Value = 1886.875
$ws_CUL.Cells.Item(74, 9).Value = 0
$ws_CUL.Cells.Item(74, 10).Value = 1886.875
$ws_CUL.Cells.Item(74, 11).Value = 0
$ws_CUL.Cells.Item(74, 12).Value = 5660.625
$ws_CUL.Cells.Item(74, 14).Value = -7782.625

# CUL row 77
$ws_CUL.Cells.Item(77, 8).Value = 1886.875
$ws_CUL.Cells.Item(77, 9).Value = 0
$ws_CUL.Cells.Item(77, 10).Value = 1886.875
$ws_CUL.Cells.Item(77, 11).Value = 0
$ws_CUL.Cells.Item(77, 12).Value = 16981.875
$ws_CUL.Cells.Item(77, 14).Value = -27589.875

# CUL row 132
$ws_CUL.Cells.Item(132, 8).Value = 1688.6428
$ws_CUL.Cells.Item(132, 9).Value = 1591.6923
$ws_CUL.Cells.Item(132, 10).Value = 1772.6666
$ws_CUL.Cells.Item(132, 11).Value = 14325.2307
$ws_CUL.Cells.Item(132, 12).Value = 15953.9994
$ws_CUL.Cells.Item(132, 13).Value = -11795.2307
$ws_CUL.Cells.Item(132, 14).Value = -21013.9994

# GSM row 32
$ws_GSM.Cells.Item(32, 8).Value = 43266.668
$ws_GSM.Cells.Item(32, 9).Value = 0
$ws_GSM.Cells.Item(32, 10).Value = 43266.668
$ws_GSM.Cells.Item(32, 11).Value = 0
$ws_GSM.Cells.Item(32, 12).Value = 43266.668
$ws_GSM.Cells.Item(32, 14).Value = -43858.668

# GSM row 97
$ws_GSM.Cells.Item(97, 8).Value = 1455.4482
$ws_GSM.Cells.Item(97, 9).Value = 1171.96
$ws_GSM.Cells.Item(97, 10).Value = 3227.25
$ws_GSM.Cells.Item(97, 11).Value = 1171.96
$ws_GSM.Cells.Item(97, 12).Value = 3227.25
$ws_GSM.Cells.Item(97, 13).Value = -675.96
$ws_GSM.Cells.Item(97, 14).Value = -4219.25

# LTW row 100
$ws_LTW.Cells.Item(100, 8).Value = 4343.3887
$ws_LTW.Cells.Item(100, 9).Value = 2598.4167
$ws_LTW.Cells.Item(100, 10).Value = 7833.3335
$ws_LTW.Cells.Item(100, 11).Value = 2598.4167
$ws_LTW.Cells.Item(100, 12).Value = 7833.3335
$ws_LTW.Cells.Item(100, 13).Value = -2057.4167
$ws_LTW.Cells.Item(100, 14).Value = -8915.333500000001

# LTW row 114
$ws_LTW.Cells.Item(114, 8).Value = 0
$ws_LTW.Cells.Item(114, 9).Value = 0
$ws_LTW.Cells.Item(114, 10).Value = 0
$ws_LTW.Cells.Item(114, 11).Value = 0
$ws_LTW.Cells.Item(114, 12).Value = 0
$ws_LTW.Cells.Item(114, 14).ClearContents()

# WVR row 63
$ws_WVR.Cells.Item(63, 8).Value = 0
$ws_WVR.Cells.Item(63, 9).Value = 0
$ws_WVR.Cells.Item(63, 10).Value = 0
$ws_WVR.Cells.Item(63, 11).Value = 0
$ws_WVR.Cells.Item(63, 12).Value = 0
$ws_WVR.Cells.Item(63, 14).ClearContents()

# WVR row 64
$ws_WVR.Cells.Item(64, 8).Value = 0
$ws_WVR.Cells.Item(64, 9).Value = 0
$ws_WVR.Cells.Item(64, 10).Value = 0
$ws_WVR.Cells.Item(64, 11).Value = 0
$ws_WVR.Cells.Item(64, 12).Value = 0
$ws_WVR.Cells.Item(64, 14).ClearContents()

# WVR row 66
$ws_WVR.Cells.Item(66, 8).Value = 0
$ws_WVR.Cells.Item(66, 9).Value = 0
$ws_WVR.Cells.Item(66, 10).Value = 0
$ws_WVR.Cells.Item(66, 11).Value = 0
$ws_WVR.Cells.Item(66, 12).Value = 0
$ws_WVR.Cells.Item(66, 14).ClearContents()

# WVR row 67
$ws_WVR.Cells.Item(67, 8).Value = 0
$ws_WVR.Cells.Item(67, 9).Value = 0
$ws_WVR.Cells.Item(67, 10).Value = 0
$ws_WVR.Cells.Item(67, 11).Value = 0
$ws_WVR.Cells.Item(67, 12).Value = 0
$ws_WVR.Cells.Item(67, 14).ClearContents()

# WVR row 70
$ws_WVR.Cells.Item(70, 8).Value = 13000
$ws_WVR.Cells.Item(70, 9).Value = 0
$ws_WVR.Cells.Item(70, 10).Value = 13000
$ws_WVR.Cells.Item(70, 11).Value = 0
$ws_WVR.Cells.Item(70, 12).Value = 13000
$ws_WVR.Cells.Item(70, 14).Value = -13630

# WVR row 73
$ws_WVR.Cells.Item(73, 8).Value = 13000
$ws_WVR.Cells.Item(73, 9).Value = 0
$ws_WVR.Cells.Item(73, 10).Value = 13000
$ws_WVR.Cells.Item(73, 11).Value = 0
$ws_WVR.Cells.Item(73, 12).Value = 13000
$ws_WVR.Cells.Item(73, 14).Value = -15184

# WVR row 76
$ws_WVR.Cells.Item(76, 8).Value = 0
$ws_WVR.Cells.Item(76, 9).Value = 0
$ws_WVR.Cells.Item(76, 10).Value = 0
$ws_WVR.Cells.Item(76, 11).Value = 0
$ws_WVR.Cells.Item(76, 12).Value = 0
$ws_WVR.Cells.Item(76, 13).ClearContents()
$ws_WVR.Cells.Item(76, 14).ClearContents()

# WVR row 79
$ws_WVR.Cells.Item(79, 8).Value = 0
$ws_WVR.Cells.Item(79, 9).Value = 0
$ws_WVR.Cells.Item(79, 10).Value = 0
$ws_WVR.Cells.Item(79, 11).Value = 0
$ws_WVR.Cells.Item(79, 12).Value = 0
$ws_WVR.Cells.Item(79, 13).ClearContents()
$ws_WVR.Cells.Item(79, 14).ClearContents()

# WVR row 80
$ws_WVR.Cells.Item(80, 8).Value = 21787
$ws_WVR.Cells.Item(80, 9).Value = 7273
$ws_WVR.Cells.Item(80, 10).Value = 36301
$ws_WVR.Cells.Item(80, 11).Value = 7273
$ws_WVR.Cells.Item(80, 12).Value = 36301
$ws_WVR.Cells.Item(80, 13).Value = -6275
$ws_WVR.Cells.Item(80, 14).Value = -38297

# WVR row 81
$ws_WVR.Cells.Item(81, 8).Value = 2972.111
$ws_WVR.Cells.Item(81, 9).Value = 2575
$ws_WVR.Cells.Item(81, 10).Value = 3289.8
$ws_WVR.Cells.Item(81, 11).Value = 5150
$ws_WVR.Cells.Item(81, 12).Value = 6579.6
$ws_WVR.Cells.Item(81, 13).Value = -4089
$ws_WVR.Cells.Item(81, 14).Value = -8701.6

# WVR row 82
$ws_WVR.Cells.Item(82, 8).Value = 0
$ws_WVR.Cells.Item(82, 9).Value = 0
$ws_WVR.Cells.Item(82, 10).Value = 0
$ws_WVR.Cells.Item(82, 11).Value = 0
$ws_WVR.Cells.Item(82, 12).Value = 0
$ws_WVR.Cells.Item(82, 14).ClearContents()

# WVR row 83
$ws_WVR.Cells.Item(83, 8).Value = 21787
$ws_WVR.Cells.Item(83, 9).Value = 7273
$ws_WVR.Cells.Item(83, 10).Value = 36301
$ws_WVR.Cells.Item(83, 11).Value = 21819
$ws_WVR.Cells.Item(83, 12).Value = 108903
$ws_WVR.Cells.Item(83, 13).Value = -16827
$ws_WVR.Cells.Item(83, 14).Value = -118887

# WVR row 84
$ws_WVR.Cells.Item(84, 8).Value = 2972.111
$ws_WVR.Cells.Item(84, 9).Value = 2575
$ws_WVR.Cells.Item(84, 10).Value = 3289.8
$ws_WVR.Cells.Item(84, 11).Value = 25750
$ws_WVR.Cells.Item(84, 12).Value = 32898
$ws_WVR.Cells.Item(84, 13).Value = -20446
$ws_WVR.Cells.Item(84, 14).Value = -43506

# WVR row 85
$ws_WVR.Cells.Item(85, 8).Value = 0
$ws_WVR.Cells.Item(85, 9).Value = 0
$ws_WVR.Cells.Item(85, 10).Value = 0
$ws_WVR.Cells.Item(85, 11).Value = 0
$ws_WVR.Cells.Item(85, 12).Value = 0
$ws_WVR.Cells.Item(85, 14).ClearContents()

# WVR row 88
$ws_WVR.Cells.Item(88, 8).Value = 30189
$ws_WVR.Cells.Item(88, 9).Value = 0
$ws_WVR.Cells.Item(88, 10).Value = 30189
$ws_WVR.Cells.Item(88, 11).Value = 0
$ws_WVR.Cells.Item(88, 12).Value = 30189
$ws_WVR.Cells.Item(88, 14).Value = -31001

# WVR row 91
$ws_WVR.Cells.Item(91, 8).Value = 30189
$ws_WVR.Cells.Item(91, 9).Value = 0
$ws_WVR.Cells.Item(91, 10).Value = 30189
$ws_WVR.Cells.Item(91, 11).Value = 0
$ws_WVR.Cells.Item(91, 12).Value = 30189
$ws_WVR.Cells.Item(91, 14).Value = -32997

# WVR row 107
$ws_WVR.Cells.Item(107, 8).Value = 1164
$ws_WVR.Cells.Item(107, 9).Value = 382.6154
$ws_WVR.Cells.Item(107, 10).Value = 2615.1428
$ws_WVR.Cells.Item(107, 11).Value = 1147.8462
$ws_WVR.Cells.Item(107, 12).Value = 7845.428400000001
$ws_WVR.Cells.Item(107, 13).Value = 772.1538
$ws_WVR.Cells.Item(107, 14).Value = -11685.4284

# WVR row 111
$ws_WVR.Cells.Item(111, 8).Value = 52632.332
$ws_WVR.Cells.Item(111, 9).Value = 0
$ws_WVR.Cells.Item(111, 10).Value = 52632.332
$ws_WVR.Cells.Item(111, 11).Value = 0
$ws_WVR.Cells.Item(111, 12).Value = 52632.332
$ws_WVR.Cells.Item(111, 14).Value = -60812.332

# WVR row 120
$ws_WVR.Cells.Item(120, 8).Value = 69800
$ws_WVR.Cells.Item(120, 9).Value = 0
$ws_WVR.Cells.Item(120, 10).Value = 69800
$ws_WVR.Cells.Item(120, 11).Value = 0
$ws_WVR.Cells.Item(120, 12).Value = 69800
$ws_WVR.Cells.Item(120, 14).Value = -79476

